$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 4
$ws.Range("I2").Value = 2.14
$ws.Range("N2").Value = 3.35
$ws.Range("R2").Value = 1.31
$ws.Range("U2").Value = 2.04
$ws.Range("V2").Value = 1.87
$ws.Range("X2").Value = 970
$ws.Range("Z2").Value = 13
$ws.Range("AC2").Value = 7.8
$ws.Range("AD2").Value = 11
$ws.Range("AE2").Value = 980
$ws.Range("AF2").Value = 980
$ws.Range("AG2").Value = 970
$ws.Range("AK2").Value = 60
$ws.Range("AO2").Value = 1000
$ws.Range("F3").Value = 1.42
$ws.Range("G3").Value = 1.65
$ws.Range("H3").Value = 5.4
$ws.Range("I3").Value = 13
$ws.Range("K3").Value = 980
$ws.Range("N3").Value = 1.89
$ws.Range("P3").Value = 1.89
$ws.Range("Q3").Value = 1.81
$ws.Range("R3").Value = 1.33
$ws.Range("S3").Value = 2.82
$ws.Range("T3").Value = 1.86
$ws.Range("V3").Value = 1.08
$ws.Range("W3").Value = 2.52
$ws.Range("Y3").Value = 1000
$ws.Range("AB3").Value = 1000
$ws.Range("AC3").Value = 1000
$ws.Range("AD3").Value = 1000
$ws.Range("AF3").Value = 1000
$ws.Range("AG3").Value = 1000
$ws.Range("AH3").Value = 1000
$ws.Range("AJ3").Value = 1000
$ws.Range("AK3").Value = 1000
$ws.Range("AN3").Value = 1000
$ws.Range("M4").Value = 1.07
$ws.Range("Q4").Value = 1.94
$ws.Range("AB4").Value = 9.4
$ws.Range("AM4").Value = 100
$ws.Range("AN4").Value = 14
$ws.Range("P5").Value = 1.99
$ws.Range("AB5").Value = 7
$ws.Range("AI5").Value = 230
$ws.Range("F6").Value = 2.46
$ws.Range("G6").Value = 3.35
$ws.Range("H6").Value = 2.5
$ws.Range("I6").Value = 3.15
$ws.Range("J6").Value = 3.3
$ws.Range("R6").Value = 1.34
$ws.Range("S6").Value = 2.8
$ws.Range("T6").Value = 1.55
$ws.Range("V6").Value = 1.46
$ws.Range("W6").Value = 1.43
$ws.Range("Y6").Value = 1000
$ws.Range("Z6").Value = 1000
$ws.Range("AB6").Value = 1000
$ws.Range("AC6").Value = 1000
$ws.Range("AD6").Value = 1000
$ws.Range("AE6").Value = 1000
$ws.Range("AF6").Value = 1000
$ws.Range("AG6").Value = 1000
$ws.Range("AH6").Value = 1000
$ws.Range("AK6").Value = 1000
$ws.Range("F7").Value = 1.53
$ws.Range("G7").Value = 1.82
$ws.Range("H7").Value = 4.7
$ws.Range("J7").Value = 4
$ws.Range("Q7").Value = 1.65
$ws.Range("V7").Value = 1.16
$ws.Range("V8").Value = 4.6
$ws.Range("AB8").Value = 80
$ws.Range("AG8").Value = 42
$ws.Range("AJ8").Value = 480
$ws.Range("G9").Value = 1.93
$ws.Range("N9").Value = 5.3
$ws.Range("U9").Value = 2.5
$ws.Range("V9").Value = 1.3
$ws.Range("X9").Value = 22
$ws.Range("AA9").Value = 80
$ws.Range("AC9").Value = 9.6
$ws.Range("AD9").Value = 16.5
$ws.Range("AG9").Value = 10.5
$ws.Range("AH9").Value = 16
$ws.Range("AI9").Value = 44
$ws.Range("AK9").Value = 17
$ws.Range("AL9").Value = 26
$ws.Range("AM9").Value = 65
$ws.Range("AN9").Value = 8.800000000000001
$ws.Range("AO9").Value = 32
$ws.Range("G10").Value = 1.8
$ws.Range("AA10").Value = 100
$ws.Range("AM10").Value = 70
$ws.Range("P11").Value = 2.42
$ws.Range("R11").Value = 1.52
$ws.Range("S11").Value = 2.26
$ws.Range("T11").Value = 1.45
$ws.Range("U11").Value = 2.28
$ws.Range("Y11").Value = 1000
$ws.Range("Z11").Value = 1000
$ws.Range("AA11").Value = 1000
$ws.Range("AB11").Value = 1000
$ws.Range("AC11").Value = 1000
$ws.Range("AD11").Value = 1000
$ws.Range("AE11").Value = 1000
$ws.Range("AF11").Value = 1000
$ws.Range("AG11").Value = 1000
$ws.Range("AH11").Value = 1000
$ws.Range("AI11").Value = 1000
$ws.Range("AJ11").Value = 1000
$ws.Range("AK11").Value = 1000
$ws.Range("AL11").Value = 1000
$ws.Range("AM11").Value = 1000
$ws.Range("AN11").Value = 1000
$ws.Range("G12").Value = 1.53
$ws.Range("N12").Value = 5.3
$ws.Range("O12").Value = 1.19
$ws.Range("R12").Value = 1.58
$ws.Range("S12").Value = 2.44
$ws.Range("T12").Value = 1.72
$ws.Range("U12").Value = 2.1
$ws.Range("W12").Value = 2.9
$ws.Range("X12").Value = 30
$ws.Range("Z12").Value = 70
$ws.Range("AA12").Value = 230
$ws.Range("AB12").Value = 11.5
$ws.Range("AC12").Value = 970
$ws.Range("AE12").Value = 110
$ws.Range("AF12").Value = 970
$ws.Range("AG12").Value = 970
$ws.Range("AH12").Value = 970
$ws.Range("AI12").Value = 80
$ws.Range("AJ12").Value = 970
$ws.Range("AK12").Value = 970
$ws.Range("AM12").Value = 110
$ws.Range("AN12").Value = 6.4
$ws.Range("AO12").Value = 110
$ws.Range("N13").Value = 6.2
$ws.Range("P13").Value = 2.7
$ws.Range("R13").Value = 1.68
$ws.Range("AG13").Value = 15
$ws.Range("AH13").Value = 65
$ws.Range("AI13").Value = 520
$ws.Range("AK13").Value = 16.5
$ws.Range("AM13").Value = 580
$ws.Range("AN13").Value = 3.4
$ws.Range("F14").Value = 1.44
$ws.Range("Y14").Value = 29
$ws.Range("AJ14").Value = 11.5
$ws.Range("S15").Value = 2.5
$ws.Range("X15").Value = 24
$ws.Range("Y15").Value = 11
$ws.Range("AA15").Value = 14
$ws.Range("AD15").Value = 10.5
$ws.Range("AE15").Value = 15
$ws.Range("AJ15").Value = 190
$ws.Range("AK15").Value = 95
$ws.Range("AL15").Value = 85
$ws.Range("AO15").Value = 6.2
$ws.Range("F16").Value = 1.26
$ws.Range("G16").Value = 1.46
$ws.Range("I16").Value = 19.5
$ws.Range("J16").Value = 5
$ws.Range("N16").Value = 1.04
$ws.Range("P16").Value = 1.45
$ws.Range("Q16").Value = 1.2
$ws.Range("R16").Value = 1.45
$ws.Range("S16").Value = 2.34
$ws.Range("T16").Value = 1.89
$ws.Range("V16").Value = 1.05
$ws.Range("W16").Value = 2.88
